$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Cells.Item(98, 8).Value = 1405.0741
$ws.Cells.Item(98, 9).Value = 1177.6
$ws.Cells.Item(98, 11).Value = 1177.6
$ws.Cells.Item(98, 13).Value = 320.4000000000001
$ws.Cells.Item(107, 8).Value = 1041.2222
$ws.Cells.Item(107, 9).Value = 679
$ws.Cells.Item(107, 11).Value = 679
$ws.Cells.Item(107, 13).Value = 1241
$ws.Cells.Item(122, 8).Value = 1405.0741
$ws.Cells.Item(122, 9).Value = 1177.6
$ws.Cells.Item(122, 11).Value = 3532.8
$ws.Cells.Item(122, 13).Value = -1082.8
$ws.Cells.Item(129, 8).Value = 878.7692
$ws.Cells.Item(129, 10).Value = 899.32355
$ws.Cells.Item(129, 12).Value = 2697.97065
$ws.Cells.Item(129, 14).Value = -12697.97065
$ws.Cells.Item(132, 8).Value = 1238.4736
$ws.Cells.Item(132, 9).Value = 1154.6
$ws.Cells.Item(132, 11).Value = 3463.8
$ws.Cells.Item(132, 13).Value = -933.7999999999997
$ws.Cells.Item(141, 8).Value = 7960
$ws.Cells.Item(141, 9).Value = 2794.6667
$ws.Cells.Item(141, 11).Value = 8384.000100000001
$ws.Cells.Item(141, 13).Value = -3204.000100000001

$ws = $wb.Worksheets("ARM")
$ws.Cells.Item(2, 8).Value = 1000.95
$ws.Cells.Item(2, 9).Value = 996.26666
$ws.Cells.Item(2, 10).Value = 1015
$ws.Cells.Item(2, 11).Value = 996.26666
$ws.Cells.Item(2, 12).Value = 1015
$ws.Cells.Item(2, 13).Value = -883.26666
$ws.Cells.Item(2, 14).Value = -1241
$ws.Cells.Item(32, 8).Value = 5204.2144
$ws.Cells.Item(32, 9).Value = 3754.2554
$ws.Cells.Item(32, 10).Value = 12776.223
$ws.Cells.Item(32, 11).Value = 3754.2554
$ws.Cells.Item(32, 12).Value = 12776.223
$ws.Cells.Item(32, 13).Value = -3467.2554
$ws.Cells.Item(32, 14).Value = -13350.223
$ws.Cells.Item(35, 8).Value = 5499.6665
$ws.Cells.Item(35, 9).Value = 5499.6665
$ws.Cells.Item(35, 11).Value = 5499.6665
$ws.Cells.Item(35, 13).Value = -5093.6665
$ws.Cells.Item(45, 8).Value = 1582.091
$ws.Cells.Item(45, 9).Value = 1052.25
$ws.Cells.Item(45, 10).Value = 1884.8572
$ws.Cells.Item(45, 11).Value = 1052.25
$ws.Cells.Item(45, 12).Value = 1884.8572
$ws.Cells.Item(45, 13).Value = -675.25
$ws.Cells.Item(45, 14).Value = -2638.8572
$ws.Cells.Item(110, 8).Value = 135.28572
$ws.Cells.Item(110, 9).Value = 141.16667
$ws.Cells.Item(110, 10).Value = 100
$ws.Cells.Item(110, 11).Value = 141.16667
$ws.Cells.Item(110, 12).Value = 100
$ws.Cells.Item(110, 13).Value = 1903.83333
$ws.Cells.Item(110, 14).Value = -4190
$ws.Cells.Item(116, 8).Value = 1000.95
$ws.Cells.Item(116, 9).Value = 996.26666
$ws.Cells.Item(116, 10).Value = 1015
$ws.Cells.Item(116, 11).Value = 996.26666
$ws.Cells.Item(116, 12).Value = 1015
$ws.Cells.Item(116, 13).Value = 1297.73334
$ws.Cells.Item(116, 14).Value = -5603
$ws.Cells.Item(122, 8).Value = 907.5769
$ws.Cells.Item(122, 9).Value = 778.13043
$ws.Cells.Item(122, 10).Value = 1900
$ws.Cells.Item(122, 11).Value = 2334.39129
$ws.Cells.Item(122, 12).Value = 5700
$ws.Cells.Item(122, 13).Value = 115.60871
$ws.Cells.Item(122, 14).Value = -10600

$ws = $wb.Worksheets("BSM")
$ws.Cells.Item(3, 8).Value = 1000.95
$ws.Cells.Item(3, 9).Value = 996.26666
$ws.Cells.Item(3, 10).Value = 1015
$ws.Cells.Item(3, 11).Value = 996.26666
$ws.Cells.Item(3, 12).Value = 1015
$ws.Cells.Item(3, 13).Value = -882.26666
$ws.Cells.Item(3, 14).Value = -1243
$ws.Cells.Item(105, 8).Value = 2318.7827
$ws.Cells.Item(105, 10).Value = 3776
$ws.Cells.Item(105, 12).Value = 3776
$ws.Cells.Item(105, 14).Value = -7270
$ws.Cells.Item(107, 8).Value = 566.82355
$ws.Cells.Item(107, 9).Value = 495.73334
$ws.Cells.Item(107, 10).Value = 1100
$ws.Cells.Item(107, 11).Value = 495.73334
$ws.Cells.Item(107, 12).Value = 1100
$ws.Cells.Item(107, 13).Value = 1424.26666
$ws.Cells.Item(107, 14).Value = -4940
$ws.Cells.Item(134, 8).Value = 5745.48
$ws.Cells.Item(134, 9).Value = 6641.0527
$ws.Cells.Item(134, 11).Value = 19923.1581
$ws.Cells.Item(134, 13).Value = -17388.1581

$ws = $wb.Worksheets("CRP")
$ws.Cells.Item(16, 8).Value = 765.8333
$ws.Cells.Item(16, 10).Value = 599
$ws.Cells.Item(16, 12).Value = 599
$ws.Cells.Item(16, 14).Value = -1173
$ws.Cells.Item(94, 8).Value = 1290.091
$ws.Cells.Item(94, 9).Value = 1327.8572
$ws.Cells.Item(94, 11).Value = 1327.8572
$ws.Cells.Item(94, 13).Value = -876.8571999999999
$ws.Cells.Item(105, 8).Value = 854.61536
$ws.Cells.Item(105, 9).Value = 777.7
$ws.Cells.Item(105, 11).Value = 777.7
$ws.Cells.Item(105, 13).Value = 969.3
$ws.Cells.Item(107, 8).Value = 946.8461
$ws.Cells.Item(107, 9).Value = 609.1667
$ws.Cells.Item(107, 10).Value = 4999
$ws.Cells.Item(107, 11).Value = 609.1667
$ws.Cells.Item(107, 12).Value = 4999
$ws.Cells.Item(107, 13).Value = 1310.8333
$ws.Cells.Item(107, 14).Value = -8839
$ws.Cells.Item(113, 8).Value = 765.8333
$ws.Cells.Item(113, 10).Value = 599
$ws.Cells.Item(113, 12).Value = 599
$ws.Cells.Item(113, 14).Value = -4939

$ws = $wb.Worksheets("CUL")
$ws.Cells.Item(131, 8).Value = 10937.4
$ws.Cells.Item(131, 10).Value = 11560.576
$ws.Cells.Item(131, 12).Value = 34681.728
$ws.Cells.Item(131, 14).Value = -44761.728

$ws = $wb.Worksheets("GSM")
$ws.Cells.Item(97, 8).Value = 1076.75
$ws.Cells.Item(97, 9).Value = 1123.4286
$ws.Cells.Item(97, 11).Value = 1123.4286
$ws.Cells.Item(97, 13).Value = -627.4286
$ws.Cells.Item(113, 8).Value = 1049.4615
$ws.Cells.Item(113, 9).Value = 960.5
$ws.Cells.Item(113, 10).Value = 1089
$ws.Cells.Item(113, 11).Value = 960.5
$ws.Cells.Item(113, 12).Value = 1089
$ws.Cells.Item(113, 13).Value = 1209.5
$ws.Cells.Item(113, 14).Value = -5429

$ws = $wb.Worksheets("LTW")
$ws.Cells.Item(7, 8).Value = 2759.1
$ws.Cells.Item(7, 9).Value = 2636.375
$ws.Cells.Item(7, 11).Value = 2636.375
$ws.Cells.Item(7, 13).Value = -2524.375
$ws.Cells.Item(61, 8).Value = 2134.238
$ws.Cells.Item(61, 9).Value = 2077.1667
$ws.Cells.Item(61, 10).Value = 2210.3333
$ws.Cells.Item(61, 11).Value = 2077.1667
$ws.Cells.Item(61, 12).Value = 2210.3333
$ws.Cells.Item(61, 13).Value = -1875.1667
$ws.Cells.Item(61, 14).Value = -2614.3333
$ws.Cells.Item(113, 8).Value = 2134.238
$ws.Cells.Item(113, 9).Value = 2077.1667
$ws.Cells.Item(113, 10).Value = 2210.3333
$ws.Cells.Item(113, 11).Value = 2077.1667
$ws.Cells.Item(113, 12).Value = 2210.3333
$ws.Cells.Item(113, 13).Value = 92.83329999999978
$ws.Cells.Item(113, 14).Value = -6550.3333
$ws.Cells.Item(122, 8).Value = 5412.609
$ws.Cells.Item(122, 9).Value = 4905.3125
$ws.Cells.Item(122, 11).Value = 14715.9375
$ws.Cells.Item(122, 13).Value = -12265.9375
$ws.Cells.Item(126, 8).Value = 2759.1
$ws.Cells.Item(126, 9).Value = 2636.375
$ws.Cells.Item(126, 11).Value = 7909.125
$ws.Cells.Item(126, 13).Value = -5439.125
$ws.Cells.Item(136, 8).Value = 5299.8887
$ws.Cells.Item(136, 9).Value = 4201.154
$ws.Cells.Item(136, 11).Value = 12603.462
$ws.Cells.Item(136, 13).Value = -10053.462

$ws = $wb.Worksheets("WVR")
$ws.Cells.Item(122, 8).Value = 36482.09
$ws.Cells.Item(122, 9).Value = 49526.5
$ws.Cells.Item(122, 11).Value = 148579.5
$ws.Cells.Item(122, 13).Value = -146129.5
$ws.Cells.Item(126, 8).Value = 4567.6665
$ws.Cells.Item(126, 9).Value = 3901.75
$ws.Cells.Item(126, 10).Value = 5899.5
$ws.Cells.Item(126, 11).Value = 11705.25
$ws.Cells.Item(126, 12).Value = 17698.5
$ws.Cells.Item(126, 13).Value = -9235.25
$ws.Cells.Item(126, 14).Value = -22638.5
$ws.Cells.Item(136, 8).Value = 2029.4
$ws.Cells.Item(136, 9).Value = 1882.1818
$ws.Cells.Item(136, 11).Value = 5646.5454
$ws.Cells.Item(136, 13).Value = -3096.5454
